# Added Squamous Intraepithelial Lesion State transitions (SUB-MOD 5).
#
# The "Clinical states" transition table on Sheet1 gains one additional
# row for the LSIL/HSIL sub-model:
#   - Row 7 (previously "HSIL -> Infection") is repurposed as the new
#     "LSIL -> Well" transition.
#   - A brand new row is inserted at row 10 for "HSIL -> Well", pushing
#     every following row (old row 10 "HSIL -> HSIL" onward) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: HSIL -> Infection  becomes  LSIL -> Well
$ws.Range("A7").Value = "LSIL"
$ws.Range("B7").Value = "Well"

# Insert a fresh row above old row 10, shifting rows 10:45 down to 11:46
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the HSIL -> Well transition
$ws.Range("A10").Value = "HSIL"
$ws.Range("B10").Value = "Well"

# Match the author's final selection/view state
$ws.Range("C10").Select()
